$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.189.52'
$ws.Range("E2").Value = '  -2.31%  '
$ws.Range("D3").Value = '1.873.62'
$ws.Range("E3").Value = '  -1.65%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.88'
$ws.Range("E5").Value = '  -1.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5066'
$ws.Range("E7").Value = '  +0.57%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3759'
$ws.Range("E8").Value = '  -1.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07166'
$ws.Range("E9").Value = '  -1.47%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8919'
$ws.Range("E10").Value = '  -1.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.78'
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("D12").Value = '1.885.21'
$ws.Range("E12").Value = '  -1.40%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07582'
$ws.Range("E13").Value = '  -0.90%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.340'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.65'
$ws.Range("E15").Value = '  -2.34%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.002'
$ws.Range("E16").Value = '  -0.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008548'
$ws.Range("E17").Value = '  -1.85%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.19'
$ws.Range("E18").Value = '  -2.89%  '
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("D20").Value = '27.227.14'
$ws.Range("E20").Value = '  -2.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.090'
$ws.Range("E21").Value = '  -1.36%  '
$ws.Range("D22").Value = '2.125.41'
$ws.Range("E22").Value = '  -2.54%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.65'
$ws.Range("E23").Value = '  -1.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.521'
$ws.Range("E24").Value = '  -0.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.90'
$ws.Range("E25").Value = '  -1.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.842'
$ws.Range("E26").Value = '  -1.94%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.05'
$ws.Range("E27").Value = '  -1.75%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.096'
$ws.Range("E28").Value = '  -5.74%  '
$ws.Range("E29").Value = '  -2.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.784'
$ws.Range("E30").Value = '  -2.48%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.728'
$ws.Range("E31").Value = '  +0.52%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09004'
$ws.Range("E32").Value = '  -0.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05156'
$ws.Range("E33").Value = '  -2.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.101'
$ws.Range("E34").Value = '  -3.50%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7494'
$ws.Range("E35").Value = '  -2.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.166'
$ws.Range("E36").Value = '  -4.41%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.565'
$ws.Range("E37").Value = '  +2.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02032'
$ws.Range("E38").Value = '  -1.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.050'
$ws.Range("E39").Value = '  +1.21%  '
$ws.Range("E40").Value = '  -1.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5379'
$ws.Range("E41").Value = '  -2.75%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.639'
$ws.Range("E42").Value = '  -3.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '115.37'
$ws.Range("E43").Value = '  +3.63%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.475'
$ws.Range("E44").Value = '  +0.36%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1483'
$ws.Range("E45").Value = '  -1.81%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4663'
$ws.Range("E46").Value = '  -2.98%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.002'
$ws.Range("E47").Value = '  +0.10%  '
$ws.Range("E48").Value = '  -4.83%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.575'
$ws.Range("E49").Value = '  -3.33%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '64.99'
$ws.Range("E50").Value = '  -3.74%  '
$ws.Range("B51").Value = 'Elrond'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '37.02'
$ws.Range("E51").Value = '  +0.28%  '
